$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Table on slide 16 switches from the custom "Table_0" style to
#    the built-in "No Style, No Grid" table style.
# ------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tbl = $tableShape.Table
$tbl.ApplyStyle("{402E4E3E-34CC-4AA4-8039-3ECBE95A1D15}")

# ------------------------------------------------------------------
# 2) The presentation's theme colour scheme (ppt/theme/theme1.xml,
#    used by the slide master) switches from the custom "Integral"
#    palette to the default Office palette.
# ------------------------------------------------------------------
$theme = $p.SlideMaster.Theme.ThemeElements.ThemeColorScheme

$theme.Item(1).RGB  = 0        # dk1      000000
$theme.Item(2).RGB  = 16777215 # lt1      FFFFFF
$theme.Item(3).RGB  = 6968388  # dk2      44546A
$theme.Item(4).RGB  = 15132391 # lt2      E7E6E6
$theme.Item(5).RGB  = 13998939 # accent1  5B9BD5
$theme.Item(6).RGB  = 3243501  # accent2  ED7D31
$theme.Item(7).RGB  = 10855845 # accent3  A5A5A5
$theme.Item(8).RGB  = 49407    # accent4  FFC000
$theme.Item(9).RGB  = 12874308 # accent5  4472C4
$theme.Item(10).RGB = 4697456  # accent6  70AD47
$theme.Item(11).RGB = 12673797 # hlink    0563C1
$theme.Item(12).RGB = 7491477  # folHlink 954F72
